$wb = $excel.ActiveWorkbook

# Update "想去人数" (column F) counts for rows 2, 3, 8, 9
# on both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 623
    $ws.Range("F3").Value = 473
    $ws.Range("F8").Value = 1199
    $ws.Range("F9").Value = 3919
}
